$d = $word.ActiveDocument

# --- Edit 1: "Duchess Marie Louise" -> "Duchess Maria Luisa" ------------
# (16 April 1845 entry about Paganini's remains being allowed into Parma)
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "Duchess Marie Louise and the Bishop of Parma",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Duchess Maria Luisa and the Bishop of Parma", 2)

# --- Edit 2: insert "(Cimitero Della Villetta) " before "almost five --
# years after his death." (3 May 1845 entry about Paganini's burial)
$r2 = $d.Content
$found2 = $r2.Find.Execute(
    "Villa Gaione, Parma almost five years after his death.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Villa Gaione, Parma (Cimitero Della Villetta) almost five years after his death.", 2)

# --- Edit 3: "10 April 2016" -> "10 June 2016" (closing date line) -----
$r3 = $d.Content
$found3 = $r3.Find.Execute(
    "10 April 2016",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "10 June 2016", 2)

Write-Host "Edit1: $found1, Edit2: $found2, Edit3: $found3"
